$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.ClearFormats()
}

Set-TextValue $ws "D2" '296.45'
Set-TextValue $ws "E2" '1.04%'
Set-TextValue $ws "D3" '42.16'
Set-TextValue $ws "E3" '3.71%'
Set-TextValue $ws "D4" '5.033'
Set-TextValue $ws "E4" '-0.26%'
Set-TextValue $ws "D5" '0.07585'
Set-TextValue $ws "E5" '2.54%'
Set-TextValue $ws "B6" 'GateToken'
Set-TextValue $ws "C6" 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
Set-TextValue $ws "D6" '4.395'
Set-TextValue $ws "E6" '2.68%'
Set-TextValue $ws "B7" 'FTXToken'
Set-TextValue $ws "C7" 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
Set-TextValue $ws "D7" '1.607'
Set-TextValue $ws "E7" '3.49%'
Set-TextValue $ws "B8" 'MXToken'
Set-TextValue $ws "C8" 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-TextValue $ws "D8" '0.9305'
Set-TextValue $ws "E8" '0.74%'
Set-TextValue $ws "B9" 'BTSEToken'
Set-TextValue $ws "C9" 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
Set-TextValue $ws "D9" '2.405'
Set-TextValue $ws "E9" '2.39%'
Set-TextValue $ws "B10" 'LiechtensteinCryptoassetsExchange'
Set-TextValue $ws "C10" 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
Set-TextValue $ws "D10" '0.1208'
Set-TextValue $ws "E10" '4.88%'
Set-TextValue $ws "B11" 'WazirX'
Set-TextValue $ws "C11" 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
Set-TextValue $ws "D11" '0.1843'
Set-TextValue $ws "E11" '6.58%'
Set-TextValue $ws "B12" 'MandalaExchangeToken'
Set-TextValue $ws "C12" 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
Set-TextValue $ws "D12" '0.09007'
Set-TextValue $ws "E12" '3.30%'
Set-TextValue $ws "B13" 'BitrueCoin'
Set-TextValue $ws "C13" 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
Set-TextValue $ws "D13" '0.03990'
Set-TextValue $ws "E13" '-4.47%'
Set-TextValue $ws "B14" 'BitMartToken'
Set-TextValue $ws "C14" 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
Set-TextValue $ws "D14" '0.1053'
Set-TextValue $ws "E14" '-0.17%'
Set-TextValue $ws "B15" 'BitForexToken'
Set-TextValue $ws "C15" 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
Set-TextValue $ws "D15" '0.001279'
Set-TextValue $ws "E15" '1.18%'
Set-TextValue $ws "B16" 'TigerCash'
Set-TextValue $ws "C16" 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
Set-TextValue $ws "D16" '0.005863'
Set-TextValue $ws "E16" '-1.41%'
Set-TextValue $ws "B17" 'LEO'
Set-TextValue $ws "C17" 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
Set-TextValue $ws "D17" '3.366'
Set-TextValue $ws "E17" '-1.44%'
Set-TextValue $ws "E18" '1.12%'
Set-TextValue $ws "D19" '7.868'
Set-TextValue $ws "E19" '2.59%'
Set-TextValue $ws "D20" '0.1420'
Set-TextValue $ws "E20" '3.02%'
Set-TextValue $ws "D21" '0.3000'
Set-TextValue $ws "E21" '4.27%'
Set-TextValue $ws "E22" '5.36%'
Set-TextValue $ws "D23" '0.001265'
Set-TextValue $ws "E23" '0.48%'
Set-TextValue $ws "D24" '0.003918'
Set-TextValue $ws "E24" '0.91%'
Set-TextValue $ws "E25" '-3.69%'
Set-TextValue $ws "E26" '0.17%'
Set-TextValue $ws "E38" '3.74%'
Set-TextValue $ws "D39" '0.05217'
Set-TextValue $ws "E39" '3.79%'
Set-TextValue $ws "D40" '0.006062'
Set-TextValue $ws "E40" '10.84%'
Set-TextValue $ws "D41" '0.007791'
Set-TextValue $ws "E41" '1.38%'
Set-TextValue $ws "D43" '0.007560'
Set-TextValue $ws "E43" '3.20%'
Set-TextValue $ws "D44" '0.007259'
Set-TextValue $ws "E44" '2.18%'
Set-TextValue $ws "D45" '0.2965'
Set-TextValue $ws "E45" '-6.20%'
Set-TextValue $ws "D46" '0.00006789'
Set-TextValue $ws "E46" '6.06%'
Set-TextValue $ws "E47" '0.24%'
Set-TextValue $ws "D48" '0.04606'
Set-TextValue $ws "E48" '172.16%'
Set-TextValue $ws "D49" '0.004202'
Set-TextValue $ws "E49" '0.01%'
Set-TextValue $ws "D50" '0.00002101'
Set-TextValue $ws "E50" '0.24%'
Set-TextValue $ws "D51" '0.0002001'
Set-TextValue $ws "E51" '0.24%'

Write-Host "Applied 94 cell updates"
